$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.122.67"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.361.80"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +5.08%  "
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("E11").Value = "  -2.09%  "
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").Value = "2.784.31"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "58.040.78"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "2.356.62"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.30%  "
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "330.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "63.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.71%  "
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.48%  "
$ws.Range("E27").Value = "  -6.72%  "
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("E30").Value = "  +1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.996"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.408"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "142.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "288.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0948"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("E44").Value = "  +2.46%  "
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0222"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("E51").Value = "  -0.16%  "
